# Update gh-pages to output generated at 456a3b4
# Refreshes "想去人数" (interest count) and "最低票价" (min price) figures
# across the 展览 (exhibitions), 演出 (shows) and 全部类型 (all types) sheets.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 37595
$wsExpo.Range("G2").Value = "已售罄"
$wsExpo.Range("F7").Value = 362
$wsExpo.Range("F12").Value = 536
$wsExpo.Range("F16").Value = 647
$wsExpo.Range("F23").Value = 2518
$wsExpo.Range("F24").Value = 998
$wsExpo.Range("F27").Value = 1157
$wsExpo.Range("F29").Value = 768
$wsExpo.Range("F30").Value = 56
$wsExpo.Range("F31").Value = 1154

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 389

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 37595
$wsAll.Range("G3").Value = "已售罄"
$wsAll.Range("F9").Value = 362
$wsAll.Range("F11").Value = 389
$wsAll.Range("F18").Value = 536
$wsAll.Range("F27").Value = 647
$wsAll.Range("F34").Value = 2518
$wsAll.Range("F35").Value = 998
$wsAll.Range("F38").Value = 1157
$wsAll.Range("F41").Value = 768
$wsAll.Range("F42").Value = 56
$wsAll.Range("F43").Value = 1154
